$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 2).Value = "Bitcoin"
$ws.Cells.Item(2, 3).Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Cells.Item(2, 4).Value = "26.551.47"
$ws.Cells.Item(2, 5).Value = "  +0.11%  "
$ws.Cells.Item(3, 2).Value = "Ethereum"
$ws.Cells.Item(3, 3).Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Cells.Item(3, 4).Value = "1.811.57"
$ws.Cells.Item(3, 5).Value = "  -0.15%  "
$ws.Cells.Item(4, 2).Value = "TetherUSD"
$ws.Cells.Item(4, 3).Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextCell 4 4 "1.004"
$ws.Cells.Item(4, 5).Value = "  -0.33%  "
$ws.Cells.Item(5, 2).Value = "USDC"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextCell 5 4 "1.003"
$ws.Cells.Item(5, 5).Value = "  -0.40%  "
$ws.Cells.Item(6, 2).Value = "BNB"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextCell 6 4 "305.98"
$ws.Cells.Item(6, 5).Value = "  -0.94%  "
$ws.Cells.Item(7, 2).Value = "XRP"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell 7 4 "0.4550"
$ws.Cells.Item(7, 5).Value = "  -0.40%  "
$ws.Cells.Item(8, 2).Value = "Cardano"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell 8 4 "0.3594"
$ws.Cells.Item(8, 5).Value = "  -2.04%  "
$ws.Cells.Item(9, 2).Value = "Dogecoin"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell 9 4 "0.07107"
$ws.Cells.Item(9, 5).Value = "  -0.31%  "
$ws.Cells.Item(10, 2).Value = "Polygon"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell 10 4 "0.8930"
$ws.Cells.Item(10, 5).Value = "  +1.40%  "
$ws.Cells.Item(11, 2).Value = "TRON"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell 11 4 "0.07708"
$ws.Cells.Item(11, 5).Value = "  -0.54%  "
$ws.Cells.Item(12, 2).Value = "Solana"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell 12 4 "19.31"
$ws.Cells.Item(12, 5).Value = "  -0.27%  "
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.782.65"
$ws.Cells.Item(13, 5).Value = "  -2.20%  "
$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell 14 4 "5.255"
$ws.Cells.Item(14, 5).Value = "  -0.83%  "
$ws.Cells.Item(15, 2).Value = "Chainlink"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell 15 4 "6.294"
$ws.Cells.Item(15, 5).Value = "  -1.22%  "
$ws.Cells.Item(16, 2).Value = "Litecoin"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell 16 4 "86.73"
$ws.Cells.Item(16, 5).Value = "  +0.04%  "
$ws.Cells.Item(17, 2).Value = "BinanceUSD"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell 17 4 "1.005"
$ws.Cells.Item(17, 5).Value = "  -0.35%  "
$ws.Cells.Item(18, 2).Value = "ShibaInu"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell 18 4 "0.000008550"
$ws.Cells.Item(18, 5).Value = "  -0.51%  "
$ws.Cells.Item(19, 2).Value = "Dai"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell 19 4 "1.003"
$ws.Cells.Item(19, 5).Value = "  -0.37%  "
$ws.Cells.Item(20, 2).Value = "WrappedBTC"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(20, 4).Value = "26.579.66"
$ws.Cells.Item(20, 5).Value = "  -0.03%  "
$ws.Cells.Item(21, 2).Value = "Avalanche"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell 21 4 "14.15"
$ws.Cells.Item(21, 5).Value = "  -0.94%  "
$ws.Cells.Item(22, 2).Value = "Uniswap"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell 22 4 "4.963"
$ws.Cells.Item(22, 5).Value = "  -1.11%  "
$ws.Cells.Item(23, 2).Value = "Cosmos"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell 23 4 "10.51"
$ws.Cells.Item(23, 5).Value = "  +0.11%  "
$ws.Cells.Item(24, 2).Value = "Toncoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell 24 4 "1.928"
$ws.Cells.Item(24, 5).Value = "  -2.70%  "
$ws.Cells.Item(25, 2).Value = "Monero"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell 25 4 "151.68"
$ws.Cells.Item(25, 5).Value = "  +0.29%  "
$ws.Cells.Item(26, 2).Value = "EthereumClassic"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell 26 4 "17.80"
$ws.Cells.Item(26, 5).Value = "  -0.81%  "
$ws.Cells.Item(27, 2).Value = "LidoDAOToken"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell 27 4 "2.016"
$ws.Cells.Item(27, 5).Value = "  -3.07%  "
$ws.Cells.Item(28, 2).Value = "BitcoinCash"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell 28 4 "112.49"
$ws.Cells.Item(28, 5).Value = "  -0.46%  "
$ws.Cells.Item(29, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell 29 4 "4.833"
$ws.Cells.Item(29, 5).Value = "  -0.63%  "
$ws.Cells.Item(30, 2).Value = "Stellar"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 30 4 "0.08726"
$ws.Cells.Item(30, 5).Value = "  +0.30%  "
$ws.Cells.Item(31, 2).Value = "HuobiToken"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell 31 4 "3.129"
$ws.Cells.Item(31, 5).Value = "  +3.38%  "
$ws.Cells.Item(32, 2).Value = "ImmutableX"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell 32 4 "0.7380"
$ws.Cells.Item(32, 5).Value = "  +1.00%  "
$ws.Cells.Item(33, 2).Value = "Filecoin"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 33 4 "4.432"
$ws.Cells.Item(33, 5).Value = "  -1.53%  "
$ws.Cells.Item(34, 2).Value = "RenderToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 34 4 "2.713"
$ws.Cells.Item(34, 5).Value = "  +1.81%  "
$ws.Cells.Item(35, 2).Value = "ARBITRUM"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 35 4 "1.110"
$ws.Cells.Item(35, 5).Value = "  -0.85%  "
$ws.Cells.Item(36, 2).Value = "TrustWalletToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell 36 4 "1.071"
$ws.Cells.Item(36, 5).Value = "  -1.27%  "
$ws.Cells.Item(37, 2).Value = "VeChain"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 37 4 "0.01934"
$ws.Cells.Item(37, 5).Value = "  -1.47%  "
$ws.Cells.Item(38, 2).Value = "MXToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell 38 4 "2.912"
$ws.Cells.Item(38, 5).Value = "  +0.71%  "
$ws.Cells.Item(39, 2).Value = "Hedera"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 39 4 "0.05070"
$ws.Cells.Item(39, 5).Value = "  -1.11%  "
$ws.Cells.Item(40, 2).Value = "TheSandbox"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell 40 4 "0.5081"
$ws.Cells.Item(40, 5).Value = "  +1.52%  "
$ws.Cells.Item(41, 2).Value = "FraxShare"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell 41 4 "6.771"
$ws.Cells.Item(41, 5).Value = "  -3.07%  "
$ws.Cells.Item(42, 2).Value = "Algorand"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell 42 4 "0.1506"
$ws.Cells.Item(42, 5).Value = "  -3.09%  "
$ws.Cells.Item(43, 2).Value = "Aptos"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell 43 4 "7.999"
$ws.Cells.Item(43, 5).Value = "  -1.99%  "
$ws.Cells.Item(44, 2).Value = "Decentraland"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell 44 4 "0.4686"
$ws.Cells.Item(44, 5).Value = "  +1.75%  "
$ws.Cells.Item(45, 2).Value = "PaxDollar"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell 45 4 "1.003"
$ws.Cells.Item(45, 5).Value = "  -0.44%  "
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell 46 4 "9.952"
$ws.Cells.Item(46, 5).Value = "  -0.08%  "
$ws.Cells.Item(47, 2).Value = "Quant"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell 47 4 "99.22"
$ws.Cells.Item(47, 5).Value = "  -2.16%  "
$ws.Cells.Item(48, 2).Value = "NEARProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell 48 4 "1.567"
$ws.Cells.Item(48, 5).Value = "  -1.29%  "
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell 49 4 "0.06002"
$ws.Cells.Item(49, 5).Value = "  +0.04%  "
$ws.Cells.Item(50, 2).Value = "Aave"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell 50 4 "63.61"
$ws.Cells.Item(50, 5).Value = "  -1.36%  "
$ws.Cells.Item(51, 2).Value = "Elrond"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextCell 51 4 "35.81"
$ws.Cells.Item(51, 5).Value = "  -1.34%  "
